$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old row 2); this shifts all subsequent rows
# up by one, which reproduces the new A/B/D (date/year) values exactly.
$ws.Rows.Item(2).Delete()

# Clear cells that should now be empty (no forecast value for that cell)
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").ClearContents()

# Set recomputed forecast values
$ws.Range("E11").Value = 2.693188401769642
$ws.Range("C12").Value = 1.785377844167058
$ws.Range("E12").Value = 2.333075171696652
$ws.Range("E13").Value = 3.947916604971446
$ws.Range("C14").Value = 5.477304442308206
$ws.Range("E14").Value = 4.052456259163839
$ws.Range("E15").Value = 4.998814576944932
$ws.Range("C16").Value = 4.666532690711245
$ws.Range("E16").Value = 3.659383764712709
$ws.Range("E17").Value = 4.673582741620552
$ws.Range("C18").Value = 5.266214435142658
$ws.Range("E18").Value = 4.181342739750682
$ws.Range("E19").Value = 4.372458986620376
$ws.Range("C20").Value = 4.811826107786477
$ws.Range("E20").Value = 4.131858242365549
$ws.Range("E21").Value = 4.927320050172312
$ws.Range("C22").Value = 6.022380124455107
$ws.Range("E22").Value = 4.474956658559948
$ws.Range("C23").Value = 5.91185619417105
$ws.Range("E23").Value = 4.365509285986957
$ws.Range("C24").Value = 5.904095356703798
$ws.Range("E24").Value = 4.348199743880454
$ws.Range("E25").Value = 4.339089271348406
$ws.Range("C26").Value = 3.109393707322261
$ws.Range("E26").Value = 3.558392386986431
$ws.Range("C27").Value = 5.114185474093769
$ws.Range("E27").Value = 5.472991335528654
$ws.Range("C28").Value = 3.799522169175473
$ws.Range("E28").Value = 2.778402897289434
$ws.Range("E29").Value = 3.243024666552685
$ws.Range("C30").Value = 3.662599762249985
$ws.Range("E30").Value = 3.993000457359908
$ws.Range("C31").Value = 2.167530781895133
$ws.Range("E31").Value = 2.573593955528963
$ws.Range("E32").Value = -5.097705497973837
$ws.Range("E33").Value = 0.2915162802050064
$ws.Range("C34").Value = -0.1964516829170981
$ws.Range("E34").Value = 3.329288211255621
$ws.Range("C35").Value = 0.5766229317536675
$ws.Range("E35").Value = 4.059584075094214
$ws.Range("E36").Value = 3.96063514023246
$ws.Range("E37").Value = 3.818597641626909
$ws.Range("C38").Value = 2.256289400228262
$ws.Range("E38").Value = 3.516729866534796
$ws.Range("C39").Value = 2.288114387968587
$ws.Range("E39").Value = 3.463553906111505
$ws.Range("E40").Value = 3.020376488332777
$ws.Range("E41").Value = 0.6985632195332103
$ws.Range("C42").Value = -1.252226393550548
$ws.Range("E42").Value = 2.967032781824974
$ws.Range("C43").Value = -2.013802094285932
$ws.Range("E43").Value = 2.374210810973465
$ws.Range("E44").Value = 2.034789645219792
$ws.Range("E45").Value = 0.3452735157291054
$ws.Range("C46").Value = -0.4399034310282546
$ws.Range("E46").Value = 2.563033601911258
$ws.Range("C47").Value = -0.5865622195987186
$ws.Range("E47").Value = 2.431929210693595
$ws.Range("E48").Value = 1.969879323458756
$ws.Range("E49").Value = 1.5902148106679
$ws.Range("C50").Value = 1.06642809951869
$ws.Range("E50").Value = 2.36642828939615
$ws.Range("C51").Value = 0.7174582534189566
$ws.Range("E51").Value = 2.061048937680932
$ws.Range("E52").Value = 1.743520202089877

Write-Host "Done applying naive forecaster bugfix edits"
